$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with default (style index 0) formatting, used to reset
# style after forcing text entry via a leading apostrophe (clears any
# quote-prefix flag picked up from the text assignment).
$normalStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = '''25.972.50'
$ws.Range("D2").Style = $normalStyle
$ws.Range("E2").Value = '  +0.76%  '

$ws.Range("D3").Value = '''1.748.13'
$ws.Range("D3").Style = $normalStyle
$ws.Range("E3").Value = '  -0.18%  '

$ws.Range("D4").Value = '''1.001'
$ws.Range("D4").Style = $normalStyle
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '''235.91'
$ws.Range("D5").Style = $normalStyle
$ws.Range("E5").Value = '  -0.44%  '

$ws.Range("D6").Value = '''1.000'
$ws.Range("D6").Style = $normalStyle
$ws.Range("E6").Value = '  +0.05%  '

$ws.Range("D7").Value = '''0.5222'
$ws.Range("D7").Style = $normalStyle
$ws.Range("E7").Value = '  +3.03%  '

$ws.Range("D8").Value = '''0.2836'
$ws.Range("D8").Style = $normalStyle
$ws.Range("E8").Value = '  +5.62%  '

$ws.Range("D9").Value = '''39.40'
$ws.Range("D9").Style = $normalStyle
$ws.Range("E9").Value = '  -3.68%  '

$ws.Range("D10").Value = '''0.06149'
$ws.Range("D10").Style = $normalStyle
$ws.Range("E10").Value = '  -0.90%  '

$ws.Range("D11").Value = '''1.754.69'
$ws.Range("D11").Style = $normalStyle
$ws.Range("E11").Value = '  +0.15%  '

$ws.Range("D12").Value = '''0.07033'
$ws.Range("D12").Style = $normalStyle
$ws.Range("E12").Value = '  +1.55%  '

$ws.Range("D13").Value = '''15.50'
$ws.Range("D13").Style = $normalStyle
$ws.Range("E13").Value = '  -0.28%  '

$ws.Range("D14").Value = '''0.6470'
$ws.Range("D14").Style = $normalStyle
$ws.Range("E14").Value = '  +4.08%  '

$ws.Range("D15").Value = '''4.541'
$ws.Range("D15").Style = $normalStyle
$ws.Range("E15").Value = '  +1.14%  '

$ws.Range("D16").Value = '''77.60'
$ws.Range("D16").Style = $normalStyle
$ws.Range("E16").Value = '  -1.32%  '

$ws.Range("D17").Value = '''1.000'
$ws.Range("D17").Style = $normalStyle
$ws.Range("E17").Value = '  +0.02%  '

$ws.Range("D18").Value = '''1.0000'
$ws.Range("D18").Style = $normalStyle
$ws.Range("E18").Value = '  +0.03%  '

$ws.Range("D19").Value = '''25.980.29'
$ws.Range("D19").Style = $normalStyle

$ws.Range("D20").Value = '''11.51'
$ws.Range("D20").Style = $normalStyle
$ws.Range("E20").Value = '  -1.54%  '

$ws.Range("D21").Value = '''0.000006637'
$ws.Range("D21").Style = $normalStyle
$ws.Range("E21").Value = '  -1.42%  '

$ws.Range("D22").Value = '''1.977.34'
$ws.Range("D22").Style = $normalStyle
$ws.Range("E22").Value = '  +0.20%  '

$ws.Range("D23").Value = '''4.174'
$ws.Range("D23").Style = $normalStyle
$ws.Range("E23").Value = '  +2.83%  '

$ws.Range("D24").Value = '''8.659'
$ws.Range("D24").Style = $normalStyle
$ws.Range("E24").Value = '  +4.97%  '

$ws.Range("D25").Value = '''5.167'
$ws.Range("D25").Style = $normalStyle
$ws.Range("E25").Value = '  -0.46%  '

$ws.Range("D26").Value = '''139.24'
$ws.Range("D26").Style = $normalStyle
$ws.Range("E26").Value = '  +1.68%  '

$ws.Range("D27").Value = '''1.508'
$ws.Range("D27").Style = $normalStyle
$ws.Range("E27").Value = '  +3.32%  '

$ws.Range("D28").Value = '''1.842'
$ws.Range("D28").Style = $normalStyle
$ws.Range("E28").Value = '  +2.39%  '

$ws.Range("D29").Value = '''15.12'
$ws.Range("D29").Style = $normalStyle
$ws.Range("E29").Value = '  -0.44%  '

$ws.Range("D30").Value = '''102.94'
$ws.Range("D30").Style = $normalStyle
$ws.Range("E30").Value = '  +0.28%  '

$ws.Range("D31").Value = '''0.08327'
$ws.Range("D31").Style = $normalStyle
$ws.Range("E31").Value = '  +0.56%  '

$ws.Range("D32").Value = '''3.673'
$ws.Range("D32").Style = $normalStyle
$ws.Range("E32").Value = '  -1.64%  '

$ws.Range("D33").Value = '''3.443'
$ws.Range("D33").Style = $normalStyle
$ws.Range("E33").Value = '  +0.32%  '

$ws.Range("D34").Value = '''0.04486'
$ws.Range("D34").Style = $normalStyle
$ws.Range("E34").Value = '  +1.53%  '

$ws.Range("E35").Value = '  -1.31%  '

$ws.Range("D36").Value = '''0.9870'
$ws.Range("D36").Style = $normalStyle
$ws.Range("E36").Value = '  -1.98%  '

$ws.Range("D37").Value = '''0.6127'
$ws.Range("D37").Style = $normalStyle
$ws.Range("E37").Value = '  +1.41%  '

$ws.Range("D38").Value = '''2.686'
$ws.Range("D38").Style = $normalStyle
$ws.Range("E38").Value = '  -0.08%  '

$ws.Range("D39").Value = '''0.01591'
$ws.Range("D39").Style = $normalStyle
$ws.Range("E39").Value = '  +1.89%  '

$ws.Range("D40").Value = '''1.947'
$ws.Range("D40").Style = $normalStyle
$ws.Range("E40").Value = '  -0.76%  '

$ws.Range("E41").Value = '  -0.06%  '

$ws.Range("D42").Value = '''101.11'
$ws.Range("D42").Style = $normalStyle
$ws.Range("E42").Value = '  -1.06%  '

$ws.Range("D43").Value = '''0.3877'
$ws.Range("D43").Style = $normalStyle
$ws.Range("E43").Value = '  +0.77%  '

$ws.Range("D44").Value = '''5.080'
$ws.Range("D44").Style = $normalStyle
$ws.Range("E44").Value = '  +4.73%  '

$ws.Range("D45").Value = '''0.7360'
$ws.Range("D45").Style = $normalStyle
$ws.Range("E45").Value = '  -2.30%  '

$ws.Range("D46").Value = '''0.05473'
$ws.Range("D46").Style = $normalStyle
$ws.Range("E46").Value = '  -0.59%  '

$ws.Range("D47").Value = '''6.319'
$ws.Range("D47").Style = $normalStyle
$ws.Range("E47").Value = '  +5.83%  '

$ws.Range("D48").Value = '''0.1119'
$ws.Range("D48").Style = $normalStyle
$ws.Range("E48").Value = '  +2.38%  '

$ws.Range("D49").Value = '''53.10'
$ws.Range("D49").Style = $normalStyle
$ws.Range("E49").Value = '  +0.89%  '

$ws.Range("D50").Value = '''30.07'
$ws.Range("D50").Style = $normalStyle
$ws.Range("E50").Value = '  -0.73%  '

$ws.Range("D51").Value = '''7.616'
$ws.Range("D51").Style = $normalStyle
$ws.Range("E51").Value = '  +1.65%  '
